$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths: A ~= 11 chars, B ~= 9.54 chars
$ws.Columns.Item(1).ColumnWidth = 10.14
$ws.Columns.Item(2).ColumnWidth = 8.6

# Row 2: a blank cell A2 pre-formatted as currency (INR) ready for a future value
$ws.Range("A2").NumberFormat = '_ "₹"\ * #,##0.00_ ;_ "₹"\ * \-#,##0.00_ ;_ "₹"\ * "-"??_ ;_ @_ '

# Row 3: a blank cell C3 pre-formatted as a date
$ws.Range("C3").NumberFormat = "yyyy\-mm\-dd"

# Row 4: new attendance record
$ws.Range("A4").Value = "KAIF"
$ws.Range("B4").Value = "22:23:24"
$ws.Range("C4").Value = 45567
$ws.Range("C4").NumberFormat = "yyyy-mm-dd"

# Restore zoom/selection to match the freshly-resaved workbook
$ws.Range("A1").Select()
$excel.ActiveWindow.Zoom = 199
